$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.351.88'
$ws.Range('E2').Value = '  -0.29%  '

$ws.Range('D3').Value = '1.908.19'
$ws.Range('E3').Value = '  +1.63%  '

$ws.Range('E4').Value = '  -0.38%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.33'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.52%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.666'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.52%  '

$ws.Range('E7').Value = '  -0.31%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.66'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.84%  '

$ws.Range('E9').Value = '  +4.30%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '53.23'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +13.20%  '

$ws.Range('E11').Value = '  +2.70%  '

$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').Value = '2.187.69'
$ws.Range('E13').Value = '  +1.74%  '

$ws.Range('E14').Value = '  +6.03%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.701'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.11%  '

$ws.Range('D16').Value = '1.891.89'
$ws.Range('E16').Value = '  +0.61%  '

$ws.Range('E17').Value = '  +1.33%  '

$ws.Range('D18').Value = '35.390.74'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.16'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.36%  '

$ws.Range('D20').Value = '0.0₃0820'
$ws.Range('E20').Value = '  +2.04%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '241.41'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.83%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.48'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.82'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.09%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.41'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +27.30%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.77%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '171.13'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.36%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.45'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.28%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.36'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.84%  '

$ws.Range('E30').Value = '  +1.58%  '

$ws.Range('D31').Value = '4.155.83'
$ws.Range('E31').Value = '  +21.72%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.14'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.32%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0566'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.07%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.944'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +12.88%  '

$ws.Range('E35').Value = '  -0.34%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.09'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.72%  '

$ws.Range('E37').Value = '  -4.00%  '

$ws.Range('E38').Value = '  -0.49%  '

$ws.Range('E39').Value = '  +1.13%  '

$ws.Range('E40').Value = '  -0.31%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0652'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +8.04%  '

$ws.Range('E42').Value = '  +1.68%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.27'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +6.50%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '89.92'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.29%  '

$ws.Range('D45').Value = '1.339.87'
$ws.Range('E45').Value = '  -1.05%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.40'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.48%  '

$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '48.84'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +38.70%  '

$ws.Range('E48').Value = '  +2.24%  '

$ws.Range('E49').Value = '  -0.93%  '

$ws.Range('E50').Value = '  -2.37%  '

$ws.Range('D51').Value = '2.096.00'
$ws.Range('E51').Value = '  +1.73%  '
